$wb = $excel.ActiveWorkbook

# ===== Sheet: ALC =====
$ws = $wb.Worksheets.Item("ALC")
# Row 6 (G6=4564)
$ws.Range("H6").Value = 92.125
$ws.Range("I6").Value = 102.4
$ws.Range("K6").Value = 307.2
$ws.Range("M6").Value = -195.2
# Row 8 (G8=4565)
$ws.Range("H8").Value = 587.125
$ws.Range("I8").Value = 587.125
$ws.Range("K8").Value = 1761.375
$ws.Range("M8").Value = -1622.375
# Row 9 (G9=5487)
$ws.Range("H9").Value = 133
$ws.Range("I9").Value = 126.6
$ws.Range("J9").Value = 149
$ws.Range("K9").Value = 126.6
$ws.Range("L9").Value = 149
$ws.Range("M9").Value = 42.40000000000001
$ws.Range("N9").Value = -487
# Row 12 (G12=5515)
$ws.Range("H12").Value = 485.15384
$ws.Range("I12").Value = 158.85715
$ws.Range("J12").Value = 865.8333
$ws.Range("K12").Value = 158.85715
$ws.Range("L12").Value = 865.8333
$ws.Range("M12").Value = 11.14285000000001
$ws.Range("N12").Value = -1205.8333
# Row 21 (G21=2149)
$ws.Range("H21").Value = 0
$ws.Range("I21").Value = 0
$ws.Range("K21").Value = 0
$ws.Range("M21").ClearContents()
# Row 23 (G23=2149)
$ws.Range("H23").Value = 0
$ws.Range("I23").Value = 0
$ws.Range("K23").Value = 0
$ws.Range("M23").ClearContents()
# Row 29 (G29=4575)
$ws.Range("H29").Value = 337.5
$ws.Range("J29").Value = 0
$ws.Range("L29").Value = 0
$ws.Range("N29").ClearContents()
# Row 34 (G34=2160)
$ws.Range("H34").Value = 10333.333
$ws.Range("J34").Value = 0
$ws.Range("L34").Value = 0
$ws.Range("N34").ClearContents()
# Row 36 (G36=2160)
$ws.Range("H36").Value = 10333.333
$ws.Range("J36").Value = 0
$ws.Range("L36").Value = 0
$ws.Range("N36").ClearContents()
# Row 38 (G38=4599)
$ws.Range("H38").Value = 839.2857
$ws.Range("I38").Value = 562.5
$ws.Range("J38").Value = 2500
$ws.Range("K38").Value = 1687.5
$ws.Range("L38").Value = 7500
$ws.Range("M38").Value = -1315.5
$ws.Range("N38").Value = -8244
# Row 58 (G58=4606)
$ws.Range("H58").Value = 3946.9092
$ws.Range("I58").Value = 127.125
$ws.Range("J58").Value = 14133
$ws.Range("K58").Value = 381.375
$ws.Range("L58").Value = 42399
$ws.Range("M58").Value = -231.375
$ws.Range("N58").Value = -42699
# Row 80 (G80=12605)
$ws.Range("H80").Value = 363.75
$ws.Range("I80").Value = 363.75
$ws.Range("K80").Value = 1091.25
$ws.Range("M80").Value = -93.25
# Row 83 (G83=12605)
$ws.Range("H83").Value = 363.75
$ws.Range("I83").Value = 363.75
$ws.Range("K83").Value = 3273.75
$ws.Range("M83").Value = 1718.25
# Row 87 (G87=10651)
$ws.Range("H87").Value = 0
$ws.Range("J87").Value = 0
$ws.Range("L87").Value = 0
$ws.Range("N87").ClearContents()
# Row 90 (G90=10651)
$ws.Range("H90").Value = 0
$ws.Range("J90").Value = 0
$ws.Range("L90").Value = 0
$ws.Range("N90").ClearContents()
# Row 96 (G96=19894)
$ws.Range("H96").Value = 819.8889
$ws.Range("I96").Value = 516.2
$ws.Range("K96").Value = 1548.6
$ws.Range("M96").Value = -175.6000000000001
# Row 97 (G97=19885)
$ws.Range("H97").Value = 7999.5835
$ws.Range("J97").Value = 8181.727
$ws.Range("L97").Value = 24545.181
$ws.Range("N97").Value = -25537.181
# Row 100 (G100=19906)
$ws.Range("H100").Value = 5639.615
$ws.Range("I100").Value = 1381.6
$ws.Range("J100").Value = 19833
$ws.Range("K100").Value = 1381.6
$ws.Range("L100").Value = 19833
$ws.Range("M100").Value = -840.5999999999999
$ws.Range("N100").Value = -20915
# Row 103 (G103=19909)
$ws.Range("H103").Value = 395
$ws.Range("I103").Value = 395
$ws.Range("K103").Value = 1185
$ws.Range("M103").Value = -599
# Row 116 (G116=27778)
$ws.Range("H116").Value = 12619
$ws.Range("I116").Value = 14856.75
$ws.Range("J116").Value = 9635.333000000001
$ws.Range("K116").Value = 14856.75
$ws.Range("L116").Value = 9635.333000000001
$ws.Range("M116").Value = -11414.75
$ws.Range("N116").Value = -16519.333
# Row 132 (G132=44049)
$ws.Range("H132").Value = 21568.139
$ws.Range("I132").Value = 26251.842
$ws.Range("J132").Value = 6847.9287
$ws.Range("K132").Value = 78755.526
$ws.Range("L132").Value = 20543.7861
$ws.Range("M132").Value = -76225.526
$ws.Range("N132").Value = -25603.7861
# Row 137 (G137=44013)
$ws.Range("H137").Value = 3089.5293
$ws.Range("I137").Value = 2320.3333
$ws.Range("K137").Value = 6960.999899999999
$ws.Range("M137").Value = -4410.999899999999
# Row 138 (G138=44169)
$ws.Range("H138").Value = 4197.8057
$ws.Range("J138").Value = 4221.5884
$ws.Range("L138").Value = 12664.7652
$ws.Range("N138").Value = -22944.7652
# Row 141 (G141=44161)
$ws.Range("H141").Value = 1658.5106
$ws.Range("I141").Value = 1607.6086
$ws.Range("J141").Value = 4000
$ws.Range("K141").Value = 4822.825800000001
$ws.Range("L141").Value = 12000
$ws.Range("M141").Value = 357.1741999999995
$ws.Range("N141").Value = -22360

# ===== Sheet: ARM =====
$ws = $wb.Worksheets.Item("ARM")
# Row 5 (G5=5091)
$ws.Range("H5").Value = 223.88889
$ws.Range("I5").Value = 243.75
$ws.Range("J5").Value = 208
$ws.Range("K5").Value = 243.75
$ws.Range("L5").Value = 208
$ws.Range("M5").Value = -131.75
$ws.Range("N5").Value = -432
# Row 32 (G32=44147)
$ws.Range("H32").Value = 2608.9124
$ws.Range("I32").Value = 2387.6428
$ws.Range("K32").Value = 2387.6428
$ws.Range("M32").Value = -2100.6428
# Row 45 (G45=27714)
$ws.Range("H45").Value = 1884.3158
$ws.Range("I45").Value = 1231.6154
$ws.Range("J45").Value = 3298.5
$ws.Range("K45").Value = 1231.6154
$ws.Range("L45").Value = 3298.5
$ws.Range("M45").Value = -854.6153999999999
$ws.Range("N45").Value = -4052.5
# Row 61 (G61=43999)
$ws.Range("H61").Value = 7982.1333
$ws.Range("I61").Value = 1394.4166
$ws.Range("K61").Value = 1394.4166
$ws.Range("M61").Value = -1182.4166
# Row 63 (G63=12528)
$ws.Range("H63").Value = 462.66666
$ws.Range("J63").Value = 0
$ws.Range("L63").Value = 0
$ws.Range("N63").ClearContents()
# Row 66 (G66=12528)
$ws.Range("H66").Value = 462.66666
$ws.Range("J66").Value = 0
$ws.Range("L66").Value = 0
$ws.Range("N66").ClearContents()
# Row 74 (G74=44000)
$ws.Range("H74").Value = 1791.8334
$ws.Range("I74").Value = 1570.174
$ws.Range("J74").Value = 2184
$ws.Range("K74").Value = 1570.174
$ws.Range("L74").Value = 2184
$ws.Range("M74").Value = -696.174
$ws.Range("N74").Value = -3932
# Row 77 (G77=44000)
$ws.Range("H77").Value = 1791.8334
$ws.Range("I77").Value = 1570.174
$ws.Range("J77").Value = 2184
$ws.Range("K77").Value = 7850.87
$ws.Range("L77").Value = 10920
$ws.Range("M77").Value = -3482.87
$ws.Range("N77").Value = -19656
# Row 88 (G88=12530)
$ws.Range("H88").Value = 21219.6
$ws.Range("I88").Value = 1600
$ws.Range("J88").Value = 26124.5
$ws.Range("K88").Value = 1600
$ws.Range("L88").Value = 26124.5
$ws.Range("M88").Value = -1194
$ws.Range("N88").Value = -26936.5
# Row 91 (G91=12530)
$ws.Range("H91").Value = 21219.6
$ws.Range("I91").Value = 1600
$ws.Range("J91").Value = 26124.5
$ws.Range("K91").Value = 1600
$ws.Range("L91").Value = 26124.5
$ws.Range("M91").Value = -196
$ws.Range("N91").Value = -28932.5
# Row 94 (G94=18055)
$ws.Range("H94").Value = 47450
$ws.Range("J94").Value = 47450
$ws.Range("L94").Value = 47450
$ws.Range("N94").Value = -49252
# Row 101 (G101=18518)
$ws.Range("H101").Value = 35000
$ws.Range("J101").Value = 35000
$ws.Range("L101").Value = 35000
$ws.Range("N101").Value = -41490
# Row 122 (G122=36168)
$ws.Range("H122").Value = 2407.963
$ws.Range("I122").Value = 1329.8889
$ws.Range("J122").Value = 4564.1113
$ws.Range("K122").Value = 3989.6667
$ws.Range("L122").Value = 13692.3339
$ws.Range("M122").Value = -1539.6667
$ws.Range("N122").Value = -18592.3339
# Row 132 (G132=43997)
$ws.Range("H132").Value = 3493.195
$ws.Range("J132").Value = 4655.8
$ws.Range("L132").Value = 13967.4
$ws.Range("N132").Value = -19027.4
# Row 136 (G136=43999)
$ws.Range("H136").Value = 7982.1333
$ws.Range("I136").Value = 1394.4166
$ws.Range("K136").Value = 4183.2498
$ws.Range("M136").Value = -1633.2498

# ===== Sheet: BSM =====
$ws = $wb.Worksheets.Item("BSM")
# Row 4 (G4=5091)
$ws.Range("H4").Value = 223.88889
$ws.Range("I4").Value = 243.75
$ws.Range("J4").Value = 208
$ws.Range("K4").Value = 243.75
$ws.Range("L4").Value = 208
$ws.Range("M4").Value = -128.75
$ws.Range("N4").Value = -438
# Row 20 (G20=14149)
$ws.Range("H20").Value = 2045.0322
$ws.Range("I20").Value = 2064.0908
$ws.Range("K20").Value = 2064.0908
$ws.Range("M20").Value = -1817.0908
# Row 64 (G64=14184)
$ws.Range("H64").Value = 840.0769
$ws.Range("I64").Value = 646.2857
$ws.Range("K64").Value = 646.2857
$ws.Range("M64").Value = -421.2857
# Row 67 (G67=14184)
$ws.Range("H67").Value = 840.0769
$ws.Range("I67").Value = 646.2857
$ws.Range("K67").Value = 646.2857
$ws.Range("M67").Value = 133.7143
# Row 82 (G82=11877)
$ws.Range("H82").Value = 51666.668
$ws.Range("I82").Value = 51666.668
$ws.Range("J82").Value = 0
$ws.Range("K82").Value = 51666.668
$ws.Range("L82").Value = 0
$ws.Range("M82").Value = -51283.668
$ws.Range("N82").ClearContents()
# Row 85 (G85=11877)
$ws.Range("H85").Value = 51666.668
$ws.Range("I85").Value = 51666.668
$ws.Range("J85").Value = 0
$ws.Range("K85").Value = 51666.668
$ws.Range("L85").Value = 0
$ws.Range("M85").Value = -50340.668
$ws.Range("N85").ClearContents()
# Row 86 (G86=12526)
$ws.Range("H86").Value = 36237
$ws.Range("I86").Value = 8510.125
$ws.Range("J86").Value = 80600
$ws.Range("K86").Value = 8510.125
$ws.Range("L86").Value = 80600
$ws.Range("M86").Value = -7387.125
$ws.Range("N86").Value = -82846
# Row 89 (G89=12526)
$ws.Range("H89").Value = 36237
$ws.Range("I89").Value = 8510.125
$ws.Range("J89").Value = 80600
$ws.Range("K89").Value = 42550.625
$ws.Range("L89").Value = 403000
$ws.Range("M89").Value = -36934.625
$ws.Range("N89").Value = -414232
# Row 105 (G105=19947)
$ws.Range("H105").Value = 4215.6665
$ws.Range("I105").Value = 3875.1667
$ws.Range("K105").Value = 3875.1667
$ws.Range("M105").Value = -2128.1667
# Row 134 (G134=43998)
$ws.Range("H134").Value = 7938.675
$ws.Range("I134").Value = 6963.1562
$ws.Range("J134").Value = 11840.75
$ws.Range("K134").Value = 20889.4686
$ws.Range("L134").Value = 35522.25
$ws.Range("M134").Value = -18354.4686
$ws.Range("N134").Value = -40592.25

# ===== Sheet: CRP =====
$ws = $wb.Worksheets.Item("CRP")
# Row 16 (G16=27691)
$ws.Range("H16").Value = 477.4
$ws.Range("I16").Value = 477.4
$ws.Range("K16").Value = 477.4
$ws.Range("M16").Value = -190.4
# Row 22 (G22=5367)
$ws.Range("H22").Value = 1158.5
$ws.Range("I22").Value = 973.6
$ws.Range("K22").Value = 973.6
$ws.Range("M22").Value = -623.6
# Row 31 (G31=44023)
$ws.Range("H31").Value = 3291.6135
$ws.Range("I31").Value = 1846.9
$ws.Range("K31").Value = 1846.9
$ws.Range("M31").Value = -1551.9
# Row 34 (G34=44023)
$ws.Range("H34").Value = 3291.6135
$ws.Range("I34").Value = 1846.9
$ws.Range("K34").Value = 1846.9
$ws.Range("M34").Value = -1644.9
# Row 58 (G58=44021)
$ws.Range("H58").Value = 929.5925999999999
$ws.Range("I58").Value = 745.7917
$ws.Range("K58").Value = 745.7917
$ws.Range("M58").Value = -542.7917
# Row 113 (G113=27691)
$ws.Range("H113").Value = 477.4
$ws.Range("I113").Value = 477.4
$ws.Range("K113").Value = 477.4
$ws.Range("M113").Value = 1692.6
# Row 132 (G132=44019)
$ws.Range("H132").Value = 3034.795
$ws.Range("I132").Value = 3256.5757
$ws.Range("J132").Value = 1815
$ws.Range("K132").Value = 9769.7271
$ws.Range("L132").Value = 5445
$ws.Range("M132").Value = -7239.7271
$ws.Range("N132").Value = -10505
# Row 134 (G134=44020)
$ws.Range("H134").Value = 2790.2104
$ws.Range("I134").Value = 2589.6572
$ws.Range("K134").Value = 7768.971600000001
$ws.Range("M134").Value = -5233.971600000001
# Row 136 (G136=44021)
$ws.Range("H136").Value = 929.5925999999999
$ws.Range("I136").Value = 745.7917
$ws.Range("K136").Value = 2237.3751
$ws.Range("M136").Value = 312.6248999999998

# ===== Sheet: CUL =====
$ws = $wb.Worksheets.Item("CUL")
# Row 2 (G2=4847)
$ws.Range("H2").Value = 110.958336
$ws.Range("I2").Value = 120.76923
$ws.Range("J2").Value = 99.36364
$ws.Range("K2").Value = 724.61538
$ws.Range("L2").Value = 596.18184
$ws.Range("M2").Value = -611.61538
$ws.Range("N2").Value = -822.18184
# Row 3 (G3=44094)
$ws.Range("H3").Value = 0
$ws.Range("I3").Value = 0
$ws.Range("K3").Value = 0
$ws.Range("M3").ClearContents()
# Row 4 (G4=4650)
$ws.Range("H4").Value = 1219218.6
$ws.Range("I4").Value = 2091893.8
$ws.Range("K4").Value = 6275681.4
$ws.Range("M4").Value = -6275569.4
# Row 5 (G5=43974)
$ws.Range("H5").Value = 505.27777
$ws.Range("I5").Value = 201.09091
$ws.Range("J5").Value = 983.2857
$ws.Range("K5").Value = 603.27273
$ws.Range("L5").Value = 2949.8571
$ws.Range("M5").Value = -491.27273
$ws.Range("N5").Value = -3173.8571
# Row 12 (G12=4854)
$ws.Range("H12").Value = 220.09091
$ws.Range("I12").Value = 72.2
$ws.Range("J12").Value = 343.33334
$ws.Range("K12").Value = 216.6
$ws.Range("L12").Value = 1030.00002
$ws.Range("M12").Value = -43.60000000000002
$ws.Range("N12").Value = -1376.00002
# Row 22 (G22=4697)
$ws.Range("H22").Value = 939
$ws.Range("J22").Value = 1066.3334
$ws.Range("L22").Value = 3199.0002
$ws.Range("N22").Value = -3537.0002
# Row 27 (G27=4697)
$ws.Range("H27").Value = 939
$ws.Range("J27").Value = 1066.3334
$ws.Range("L27").Value = 3199.0002
$ws.Range("N27").Value = -3403.0002
# Row 34 (G34=4749)
$ws.Range("H34").Value = 3512.5
$ws.Range("I34").Value = 371.33334
$ws.Range("J34").Value = 8224.25
$ws.Range("K34").Value = 1114.00002
$ws.Range("L34").Value = 24672.75
$ws.Range("M34").Value = -1030.00002
$ws.Range("N34").Value = -24840.75
# Row 39 (G39=4712)
$ws.Range("H39").Value = 5210.6
$ws.Range("I39").Value = 1950
$ws.Range("J39").Value = 5494.1304
$ws.Range("K39").Value = 5850
$ws.Range("L39").Value = 16482.3912
$ws.Range("M39").Value = -5556
$ws.Range("N39").Value = -17070.3912
# Row 52 (G52=31902)
$ws.Range("H52").Value = 250
$ws.Range("J52").Value = 250
$ws.Range("L52").Value = 750
$ws.Range("N52").Value = -1282
# Row 55 (G55=4733)
$ws.Range("H55").Value = 3540.4375
$ws.Range("I55").Value = 749
$ws.Range("K55").Value = 2247
$ws.Range("M55").Value = -2070
# Row 59 (G59=4694)
$ws.Range("H59").Value = 4333
# Row 62 (G62=12845)
$ws.Range("H62").Value = 8873.625
$ws.Range("I62").Value = 6799.75
$ws.Range("J62").Value = 10947.5
$ws.Range("K62").Value = 20399.25
$ws.Range("L62").Value = 32842.5
$ws.Range("M62").Value = -19713.25
$ws.Range("N62").Value = -34214.5
# Row 65 (G65=12845)
$ws.Range("H65").Value = 8873.625
$ws.Range("I65").Value = 6799.75
$ws.Range("J65").Value = 10947.5
$ws.Range("K65").Value = 61197.75
$ws.Range("L65").Value = 98527.5
$ws.Range("M65").Value = -57765.75
$ws.Range("N65").Value = -105391.5
# Row 81 (G81=12843)
$ws.Range("H81").Value = 7500
$ws.Range("I81").Value = 0
$ws.Range("K81").Value = 0
$ws.Range("M81").ClearContents()
# Row 82 (G82=12856)
$ws.Range("H82").Value = 12843.333
$ws.Range("J82").Value = 14270
$ws.Range("L82").Value = 42810
$ws.Range("N82").Value = -43622
# Row 84 (G84=12843)
$ws.Range("H84").Value = 7500
$ws.Range("I84").Value = 0
$ws.Range("K84").Value = 0
$ws.Range("M84").ClearContents()
# Row 85 (G85=12856)
$ws.Range("H85").Value = 12843.333
$ws.Range("J85").Value = 14270
$ws.Range("L85").Value = 42810
$ws.Range("N85").Value = -45618
# Row 109 (G109=27854)
$ws.Range("H109").Value = 2373.7144
$ws.Range("J109").Value = 0
$ws.Range("L109").Value = 0
$ws.Range("N109").ClearContents()
# Row 112 (G112=27855)
$ws.Range("H112").Value = 18888.8
$ws.Range("I112").Value = 14444
$ws.Range("K112").Value = 43332
$ws.Range("M112").Value = -42224
# Row 129 (G129=36054)
$ws.Range("H129").Value = 10419344
$ws.Range("I129").Value = 1220.75
$ws.Range("J129").Value = 13892051
$ws.Range("K129").Value = 3662.25
$ws.Range("L129").Value = 41676153
$ws.Range("M129").Value = 1337.75
$ws.Range("N129").Value = -41686153
# Row 131 (G131=36060)
$ws.Range("H131").Value = 11115743
$ws.Range("I131").Value = 55558056
$ws.Range("J131").Value = 8776674
$ws.Range("K131").Value = 166674168
$ws.Range("L131").Value = 26330022
$ws.Range("M131").Value = -166669128
$ws.Range("N131").Value = -26340102
# Row 132 (G132=43972)
$ws.Range("H132").Value = 2352.5557
$ws.Range("I132").Value = 2081.8572
$ws.Range("J132").Value = 3300
$ws.Range("K132").Value = 18736.7148
$ws.Range("L132").Value = 29700
$ws.Range("M132").Value = -16206.7148
$ws.Range("N132").Value = -34760
# Row 133 (G133=44073)
$ws.Range("H133").Value = 16050.7
$ws.Range("J133").Value = 17388.611
$ws.Range("L133").Value = 52165.833
$ws.Range("N133").Value = -62285.833
# Row 135 (G135=43974)
$ws.Range("H135").Value = 505.27777
$ws.Range("I135").Value = 201.09091
$ws.Range("J135").Value = 983.2857
$ws.Range("K135").Value = 1809.81819
$ws.Range("L135").Value = 8849.5713
$ws.Range("M135").Value = 725.18181
$ws.Range("N135").Value = -13919.5713
# Row 139 (G139=44102)
$ws.Range("H139").Value = 13898168
$ws.Range("I139").Value = 18523836
$ws.Range("J139").Value = 21166.666
$ws.Range("K139").Value = 55571508
$ws.Range("L139").Value = 63499.99800000001
$ws.Range("M139").Value = -55566368
$ws.Range("N139").Value = -73779.99800000001
# Row 140 (G140=44097)
$ws.Range("H140").Value = 4472724.5
$ws.Range("I140").Value = 15628734
$ws.Range("K140").Value = 46886202
$ws.Range("M140").Value = -46881022
# Row 141 (G141=44076)
$ws.Range("H141").Value = 22733.158
$ws.Range("I141").Value = 4967.7827
$ws.Range("J141").Value = 42190.477
$ws.Range("K141").Value = 14903.3481
$ws.Range("L141").Value = 126571.431
$ws.Range("M141").Value = -9723.348099999999
$ws.Range("N141").Value = -136931.431

# ===== Sheet: GSM =====
$ws = $wb.Worksheets.Item("GSM")
# Row 9 (G9=1683)
$ws.Range("H9").Value = 200
$ws.Range("J9").Value = 200
$ws.Range("L9").Value = 200
$ws.Range("N9").Value = -540
# Row 70 (G70=14146)
$ws.Range("H70").Value = 7587.231
$ws.Range("I70").Value = 7610.125
$ws.Range("J70").Value = 7550.6
$ws.Range("K70").Value = 7610.125
$ws.Range("L70").Value = 7550.6
$ws.Range("M70").Value = -7340.125
$ws.Range("N70").Value = -8090.6
# Row 73 (G73=14146)
$ws.Range("H73").Value = 7587.231
$ws.Range("I73").Value = 7610.125
$ws.Range("J73").Value = 7550.6
$ws.Range("K73").Value = 7610.125
$ws.Range("L73").Value = 7550.6
$ws.Range("M73").Value = -6674.125
$ws.Range("N73").Value = -9422.6
# Row 102 (G102=36169)
$ws.Range("H102").Value = 3712.1035
$ws.Range("I102").Value = 3237.1667
$ws.Range("K102").Value = 3237.1667
$ws.Range("M102").Value = -1615.1667
# Row 104 (G104=18666)
$ws.Range("H104").Value = 0
$ws.Range("J104").Value = 0
$ws.Range("L104").Value = 0
$ws.Range("N104").ClearContents()
# Row 122 (G122=36182)
$ws.Range("H122").Value = 2335.0356
$ws.Range("I122").Value = 2242.9092
$ws.Range("K122").Value = 6728.7276
$ws.Range("M122").Value = -4278.7276
# Row 132 (G132=44008)
$ws.Range("H132").Value = 2247.5
$ws.Range("I132").Value = 1961.4
$ws.Range("K132").Value = 5884.200000000001
$ws.Range("M132").Value = -3354.200000000001
# Row 136 (G136=42218)
$ws.Range("H136").Value = 36470.05
$ws.Range("J136").Value = 36470.05
$ws.Range("L136").Value = 109410.15
$ws.Range("N136").Value = -114510.15

# ===== Sheet: LTW =====
$ws = $wb.Worksheets.Item("LTW")
# Row 16 (G16=5289)
$ws.Range("H16").Value = 1324.7916
$ws.Range("I16").Value = 1301.762
$ws.Range("K16").Value = 1301.762
$ws.Range("M16").Value = -1131.762
# Row 22 (G22=5277)
$ws.Range("H22").Value = 1623.625
$ws.Range("I22").Value = 877.6
$ws.Range("K22").Value = 877.6
$ws.Range("M22").Value = -582.6
# Row 27 (G27=5277)
$ws.Range("H27").Value = 1623.625
$ws.Range("I27").Value = 877.6
$ws.Range("K27").Value = 877.6
$ws.Range("M27").Value = -770.6
# Row 30 (G30=1688)
$ws.Range("H30").Value = 1758
$ws.Range("J30").Value = 3500
$ws.Range("L30").Value = 3500
$ws.Range("N30").Value = -3716
# Row 40 (G40=36248)
$ws.Range("H40").Value = 2302.4783
$ws.Range("J40").Value = 2772
$ws.Range("L40").Value = 2772
$ws.Range("N40").Value = -3044
# Row 46 (G46=5282)
$ws.Range("H46").Value = 2899.2144
$ws.Range("I46").Value = 1432
$ws.Range("J46").Value = 3299.3635
$ws.Range("K46").Value = 1432
$ws.Range("L46").Value = 3299.3635
$ws.Range("M46").Value = -1244
$ws.Range("N46").Value = -3675.3635
# Row 55 (G55=5284)
$ws.Range("H55").Value = 493.73685
$ws.Range("I55").Value = 405.8125
$ws.Range("K55").Value = 405.8125
$ws.Range("M55").Value = -232.8125
# Row 122 (G122=36247)
$ws.Range("H122").Value = 3505.4
$ws.Range("I122").Value = 1869.2727
$ws.Range("K122").Value = 5607.8181
$ws.Range("M122").Value = -3157.8181
# Row 132 (G132=44058)
$ws.Range("H132").Value = 2849.6538
$ws.Range("I132").Value = 2555.3
$ws.Range("K132").Value = 7665.900000000001
$ws.Range("M132").Value = -5135.900000000001
# Row 136 (G136=44060)
$ws.Range("H136").Value = 1645.2778
$ws.Range("I136").Value = 1188.3846
$ws.Range("K136").Value = 3565.1538
$ws.Range("M136").Value = -1015.1538

# ===== Sheet: WVR =====
$ws = $wb.Worksheets.Item("WVR")
# Row 2 (G2=3307)
$ws.Range("H2").Value = 51891.668
$ws.Range("J2").Value = 50880
$ws.Range("L2").Value = 50880
$ws.Range("N2").Value = -51104
# Row 81 (G81=12596)
$ws.Range("H81").Value = 5715.4375
$ws.Range("I81").Value = 6111.385
$ws.Range("K81").Value = 12222.77
$ws.Range("M81").Value = -11161.77
# Row 84 (G84=12596)
$ws.Range("H84").Value = 5715.4375
$ws.Range("I84").Value = 6111.385
$ws.Range("K84").Value = 61113.85000000001
$ws.Range("M84").Value = -55809.85000000001
# Row 121 (G121=26316)
$ws.Range("H121").Value = 68332.664
$ws.Range("J121").Value = 68332.664
$ws.Range("L121").Value = 68332.664
$ws.Range("N121").Value = -71826.664
# Row 126 (G126=36210)
$ws.Range("H126").Value = 2244.8
$ws.Range("I126").Value = 1959.3334
$ws.Range("K126").Value = 5878.0002
$ws.Range("M126").Value = -3408.0002
# Row 132 (G132=44029)
$ws.Range("H132").Value = 2210.849
$ws.Range("I132").Value = 2198.0435
$ws.Range("K132").Value = 6594.130500000001
$ws.Range("M132").Value = -4064.130500000001
# Row 136 (G136=44031)
$ws.Range("H136").Value = 2922.0588
$ws.Range("I136").Value = 2769.2856
$ws.Range("K136").Value = 8307.856800000001
$ws.Range("M136").Value = -5757.856800000001
# Row 137 (G137=42184)
$ws.Range("H137").Value = 58333.332
$ws.Range("J137").Value = 58333.332
$ws.Range("L137").Value = 58333.332
$ws.Range("N137").Value = -68533.33199999999
